# Rework the "Input" worksheet (sheet 1) to match the new standard template
# layout: new header labels, reordered/merged columns, and new sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Wipe out the old A1:Q5 block entirely (contents + formatting), since the
# new layout uses fewer columns (A:P) and a different header style (no bold).
$ws.Range("A1:Q5").Clear()

# New header row (A1:P1) - plain/default style, no bold.
$headers = @("발주일자", "납기일자", "거래처명", "거래처 이메일", "납품처명", "납품처 이메일", "프로젝트명", "대분류", "중분류", "소분류", "품목명", "규격", "수량", "단가", "총금액", "비고")

for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# New data rows (A2:P5).
# Columns (1-based): 1 발주일자 2 납기일자 3 거래처명 4 거래처 이메일 5 납품처명
# 6 납품처 이메일 7 프로젝트명 8 대분류 9 중분류 10 소분류 11 품목명 12 규격
# 13 수량(number) 14 단가(number) 15 총금액(number) 16 비고
$dataRows = @(@("2025-09-06", "2025-10-05", "영세엔지텍", "영세엔지텍@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "1. 원자재비", "4) ALUM. 창호", "C. 단열", "O-035 L:5660", "KS규격-1", 164, 1400, 252560, "29EA"), @("2025-09-15", "2025-09-14", "영세엔지텍", "영세엔지텍@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "1. 원자재비", "4) ALUM. 창호", "C. 단열", "YJG-176", "KS규격-2", 160, 230, 40480, "160EA"), @("2025-09-18", "2025-10-18", "영세엔지텍", "영세엔지텍@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "1. 원자재비", "4) ALUM. 창호", "C. 단열", "YJG-694", "KS규격-3", 50, 220, 12100, "50EA"), @("2025-08-31", "2025-09-18", "영세엔지텍", "영세엔지텍@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "1. 원자재비", "4) ALUM. 창호", "C. 단열", "PA-1255", "KS규격-4", 184, 400, 80960, "1차 20EA"))

$numericCols = @(13, 14, 15)
# Only columns A (발주일자) and B (납기일자) hold date-like text ("2025-09-06")
# that Excel would otherwise silently reinterpret as a date serial number.
$dateLikeCols = @(1, 2)

for ($r = 0; $r -lt $dataRows.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $dataRows[$r]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $cell = $ws.Cells.Item($rowNum, $col)
        if ($numericCols -contains $col) {
            $cell.Value = $rowData[$col - 1]
        } else {
            if ($dateLikeCols -contains $col) {
                # Force text so date-like strings (e.g. 2025-09-06) are not
                # reinterpreted as Excel date serial numbers.
                $cell.NumberFormat = "@"
            }
            $cell.Value = $rowData[$col - 1]
        }
    }
}
